# The commit swaps the storage order of the two SharePoint-related
# customXml parts in the package: the "FormTemplates" part (content-type
# forms) and the "p:properties" part (Unified Compliance Policy /
# documentManagement) change slots, so the properties part becomes
# customXml/item1.xml and the forms part becomes customXml/item2.xml.
#
# Reproduce that with the standard CustomXMLParts COM surface: pull the
# two parts' XML out (matched by namespace, so this is robust to
# whichever slot each currently occupies), delete them, then add them
# back in the opposite order so the part that used to be added first is
# now added second (and therefore ends up renumbered the other way).

$p = $ppt.ActivePresentation
$parts = $p.CustomXMLParts

$formsNamespace = "http://schemas.microsoft.com/sharepoint/v3/contenttype/forms"
$propertiesNamespace = "http://schemas.microsoft.com/office/2006/metadata/properties"

$formsXml = $null
$propertiesXml = $null
$formsPart = $null
$propertiesPart = $null

for ($i = 1; $i -le $parts.Count; $i++) {
    $part = $parts.Item($i)
    if ($part.NamespaceURI -eq $formsNamespace) {
        $formsXml = $part.XML
        $formsPart = $part
    } elseif ($part.NamespaceURI -eq $propertiesNamespace) {
        $propertiesXml = $part.XML
        $propertiesPart = $part
    }
}

if ($formsPart -ne $null) { $formsPart.Delete() }
if ($propertiesPart -ne $null) { $propertiesPart.Delete() }

# Re-add in the swapped order: properties part first (-> item1.xml),
# forms part second (-> item2.xml).
if ($propertiesXml -ne $null) { $parts.Add($propertiesXml) | Out-Null }
if ($formsXml -ne $null) { $parts.Add($formsXml) | Out-Null }
